$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.107.09"
$ws.Range("E2").Value = "  +4.16%  "

$ws.Range("D3").Value = "4.022.66"
$ws.Range("E3").Value = "  +3.50%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.03"
$ws.Range("E5").Value = "  -1.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.17"
$ws.Range("E6").Value = "  +2.89%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.687"
$ws.Range("E7").Value = "  +12.54%  "

$ws.Range("E8").Value = "  +0.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.757"
$ws.Range("E9").Value = "  +5.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.175"
$ws.Range("E10").Value = "  +1.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000327"
$ws.Range("E11").Value = "  -1.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.81"
$ws.Range("E12").Value = "  +11.95%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.86"
$ws.Range("E13").Value = "  +6.76%  "

$ws.Range("D14").Value = "4.675.53"
$ws.Range("E14").Value = "  +3.64%  "

$ws.Range("D15").Value = "4.031.56"
$ws.Range("E15").Value = "  +3.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.22"
$ws.Range("E16").Value = "  +8.33%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.17"
$ws.Range("E17").Value = "  +1.88%  "

$ws.Range("E18").Value = "  -0.60%  "

$ws.Range("E19").Value = "  -1.70%  "

$ws.Range("D20").Value = "72.099.46"
$ws.Range("E20").Value = "  +4.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "436.69"
$ws.Range("E21").Value = "  +3.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "101.55"
$ws.Range("E22").Value = "  +15.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.55"
$ws.Range("E23").Value = "  +6.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.74"
$ws.Range("E24").Value = "  +4.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.99"
$ws.Range("E25").Value = "  +0.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.72"
$ws.Range("E26").Value = "  +2.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.18"
$ws.Range("E27").Value = "  +7.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.49"
$ws.Range("E28").Value = "  +3.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.06"
$ws.Range("E29").Value = "  +9.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.52"
$ws.Range("E30").Value = "  +3.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "685.71"
$ws.Range("E31").Value = "  -1.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.128"
$ws.Range("E32").Value = "  +2.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.90"
$ws.Range("E33").Value = "  +17.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "68.00"
$ws.Range("E34").Value = "  +0.90%  "

$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.61"
$ws.Range("E35").Value = "  +4.77%  "

$ws.Range("B36").Value = "TheGraph"
$ws.Range("C36").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.439"
$ws.Range("E36").Value = "  +1.29%  "

$ws.Range("D37").Value = "0.0₃0875"
$ws.Range("E37").Value = "  +5.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.55"
$ws.Range("E38").Value = "  +19.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.152"
$ws.Range("E39").Value = "  +1.77%  "

$ws.Range("E40").Value = "  -0.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0490"
$ws.Range("E42").Value = "  +2.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.25"
$ws.Range("E43").Value = "  +5.92%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.156"
$ws.Range("E44").Value = "  +11.85%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.75"
$ws.Range("E45").Value = "  +0.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.51"
$ws.Range("E46").Value = "  +5.74%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.08"
$ws.Range("E47").Value = "  +2.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.03"
$ws.Range("E48").Value = "  +8.93%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000269"
$ws.Range("E49").Value = "  +21.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.29"
$ws.Range("E50").Value = "  +0.77%  "

$ws.Range("D51").Value = "0.0₆0339"
$ws.Range("E51").Value = "  -0.51%  "
